# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    3  = 3107
    5  = 2763
    6  = 189
    7  = 144
    9  = 1470
    12 = 21
    13 = 1229
    15 = 379
    17 = 48
    18 = 43
    20 = 78
    22 = 2699
    24 = 324
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
